# Optimized the shit out of the the telemetry module
#
# sheet2 ("http-size") previously hosted two Excel Tables (Table1 @ B5:D17,
# Table13 @ G5:I17) under merged title cells in row 4. The rework rips out
# both tables, drops the now-unused "Individual POST / Array POST" header
# block, and replaces it with a compact three-column (Individual / Array /
# Compound) summary living in B6:I10, reusing the existing bold-header
# style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- 1. Remove the two ListObjects (Table1 / Table13) entirely -----------
for ($i = $ws.ListObjects.Count; $i -ge 1; $i--) {
    $ws.ListObjects.Item($i).Delete()
}

# --- 2. Drop the merged title cells above the old tables ------------------
$ws.Range("B4:D4").UnMerge()
$ws.Range("G4:I4").UnMerge()

# --- 3. Wipe out all old content/rows (3-17) cleanly, then re-insert 5
#        blank rows so the new layout starts at row 6, matching the target
#        workbook; this avoids leftover row-height / style residue that
#        Cells.Clear() alone would leave behind. ---------------------------
$ws.Range("A3:A17").EntireRow.Delete()
$ws.Range("A1:A5").EntireRow.Insert()

# --- 4. Column widths for the new layout -----------------------------------
$ws.Columns.Item(4).ColumnWidth = 8.5      # D ~ 9.29 chars
$ws.Columns.Item(5).ColumnWidth = 13.67    # E ~ 14.57 chars
$ws.Columns.Item(7).ColumnWidth = 8.33     # G ~ 9.14 chars
$ws.Columns.Item(9).ColumnWidth = 9.17     # I = 10 chars
$ws.Columns.Item(20).ColumnWidth = 13.0    # T ~ 13.86 chars

# --- 5. New header row (row 6) — bold, no border ---------------------------
$ws.Range("B6").Value = "Individual:"
$ws.Range("E6").Value = "Array"
$ws.Range("H6").Value = "Compound"
$ws.Range("B6").Font.Bold = $true
$ws.Range("E6").Font.Bold = $true
$ws.Range("H6").Font.Bold = $true

# --- 6. Individual block (B7:C9) -------------------------------------------
$ws.Range("B7").Value = "POSTs"
$ws.Range("C7").Value = 6
$ws.Range("B8").Value = "Avg Size"
$ws.Range("C8").Value = 88
$ws.Range("B9").Value = "Total Size"
$ws.Range("C9").Formula = "=C7*C8"

# --- 7. Array block (E7:F10) -----------------------------------------------
$ws.Range("E7").Value = "POSTs"
$ws.Range("F7").Value = 1
$ws.Range("E8").Value = "Elements"
$ws.Range("F8").Value = 6
$ws.Range("E9").Value = "Size per Array"
$ws.Range("F9").Value = 534
$ws.Range("E10").Value = "Size per Element"
$ws.Range("F10").Formula = "=F9/F8"

# --- 8. Compound block (H7:I9) ----------------------------------------------
$ws.Range("H7").Value = "POSTs"
$ws.Range("I7").Value = 1
$ws.Range("H8").Value = "Elements (Inserts)"
$ws.Range("I8").Value = 6
$ws.Range("H9").Value = "Size"
$ws.Range("I9").Value = 137

# --- 9. Selection matches the saved view -----------------------------------
$ws.Range("G11").Select()
